# The underlying dataset rows (2-25) were re-ordered (the weekly
# consolidation re-ran and rows landed in a new arrangement) while the
# header row (row 1) and the set of records themselves are unchanged.
# We snapshot every data row first, then write each snapshot back to its
# new destination row, so no value is overwritten before it's captured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 25
$lastCol = 20   # columns A..T

# 1) Snapshot all existing data rows (row number -> array of 20 values).
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,$ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Destination row -> source row mapping describing where each record
#    ends up after the re-consolidation.
$mapping = @{
    2  = 18
    3  = 19
    4  = 9
    5  = 5
    6  = 6
    7  = 2
    8  = 3
    9  = 14
    10 = 15
    11 = 21
    12 = 22
    13 = 7
    14 = 8
    15 = 24
    16 = 25
    17 = 10
    18 = 4
    19 = 17
    20 = 23
    21 = 20
    22 = 11
    23 = 12
    24 = 16
    25 = 13
}

# 3) Write each destination row from its captured source-row snapshot.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - 1]
    }
}
